# database/industries/darou/delor/product/yearly.xlsx
# "update database and change read_price algorithm"
#
# The yearly columns E:I hold five consecutive "twelve months ended" periods.
# A new fiscal year (1401/12) was added and the oldest one (1396/12) dropped,
# so every data block shifts one column to the left (E<-F, F<-G, G<-H, H<-I)
# and the freed-up column I receives the new year's figure. The header rows
# that label the five year columns shift the same way. Two "sale rate" unit
# labels (for "Sales returns" / "Other & discounts" rows) were also corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header rows: each has the 5 period captions in columns E..I (col 5..9).
#    Shift the captions left by one year and add the new one in column I.
# ---------------------------------------------------------------------------
$years = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
$headerRows = @(8, 38, 71, 103, 134, 147)
foreach ($hr in $headerRows) {
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item($hr, 5 + $i).Value2 = $years[$i]
    }
}

# ---------------------------------------------------------------------------
# 2) Data rows: shift E:H to hold the former F:I values, then drop in the
#    brand-new right-most (1401/12) figure for column I.
# ---------------------------------------------------------------------------
$newI = @{
    10 = "-";         11 = "-";         12 = "-";         13 = 13964953;    14 = "-";
    15 = "-";          16 = 12372937;    17 = 921932564;   18 = "-";         19 = "-";
    20 = "-";          21 = "-";         22 = "-";         23 = "-";         24 = "-";
    25 = "-";          26 = "-";         27 = "-";         28 = "-";         29 = 0;
    30 = "-";          31 = "-";         32 = 149845492;   33 = 79351689;    34 = 1177467635;
    40 = "-";          41 = "-";         42 = "-";         43 = 16942621;    44 = "-";
    45 = "-";          46 = 14107616;    47 = "-";         48 = "-";         49 = 937446977;
    50 = "-";          51 = "-";         52 = "-";         53 = "-";         54 = "-";
    55 = "-";          56 = "-";         57 = "-";         58 = "-";         59 = "-";
    60 = "-";          61 = "-";         62 = 0;           63 = "-";         64 = "-";
    65 = 150439729;    66 = 74346276;    67 = 1193283219;
    73 = "-";          74 = "-";         75 = "-";         76 = 2553317;     77 = "-";
    78 = "-";          79 = 378659;      80 = "-";         81 = 5007868;     82 = "-";
    83 = "-";          84 = "-";         85 = "-";         86 = "-";         87 = "-";
    88 = "-";          89 = "-";         90 = "-";         91 = "-";         92 = "-";
    93 = "-";          94 = 0;           95 = "-";         96 = "-";         97 = 12554855;
    98 = 4826006;      99 = 25320705;
    105 = "-";         106 = "-";        107 = "-";        108 = 150704;     109 = "-";
    110 = "-";         111 = 26841;      112 = "-";        113 = 5342;       114 = "-";
    115 = "-";         116 = "-";        117 = "-";        118 = "-";        119 = "-";
    120 = "-";         121 = "-";        122 = "-";        123 = "-";        124 = "-";
    125 = "-";         126 = 0;          127 = "-";        128 = "-";        129 = 83454;
    130 = 64913;
    136 = -1382709;    137 = -116916;    138 = "-";        139 = -2972916;   140 = 0;
    141 = -6192342;    142 = -2744565;   143 = -13409448;
    149 = 1170608;     150 = 261743;     151 = "-";        152 = 2034952;    153 = 0;
    154 = 6362513;     155 = 2081441;    156 = 11911257;
}

foreach ($r in $newI.Keys) {
    # Read the current F,G,H,I values before overwriting anything in the row.
    $f = $ws.Cells.Item($r, 6).Value2
    $g = $ws.Cells.Item($r, 7).Value2
    $h = $ws.Cells.Item($r, 8).Value2
    $i = $ws.Cells.Item($r, 9).Value2

    $ws.Cells.Item($r, 5).Value2 = $f
    $ws.Cells.Item($r, 6).Value2 = $g
    $ws.Cells.Item($r, 7).Value2 = $h
    $ws.Cells.Item($r, 8).Value2 = $i
    $ws.Cells.Item($r, 9).Value2 = $newI[$r]
}

# ---------------------------------------------------------------------------
# 3) Fix the "sale rate" unit labels for the two rows whose quantity is not
#    counted in pieces (Sales returns / Other & discounts).
# ---------------------------------------------------------------------------
$ws.Cells.Item(112, 3).Value2 = "ریال / ریال"
$ws.Cells.Item(126, 3).Value2 = "/ ریال"
